$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Activecases"

$ws.Range("A2").Value = "21st Birthday Party 27 Nov Middels Drouin"
$ws.Range("B2").Value = 10
$ws.Range("A3").Value = "3398 BlueCross Elly Kay Mordialloc"
$ws.Range("B3").Value = 41
$ws.Range("A4").Value = "3601 Baptcare Westhaven community outbreak"
$ws.Range("B4").Value = 15
$ws.Range("A5").Value = "3646 Mornington Bay Care Community MountMartha"
$ws.Range("B5").Value = 11
$ws.Range("A6").Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Range("B6").Value = 22
$ws.Range("A7").Value = "3975 Aurrum Aged Care Brunswick West"
$ws.Range("B7").Value = 12
$ws.Range("A8").Value = "4257 BlueCross The Gables Camberwell"
$ws.Range("B8").Value = 21
$ws.Range("A9").Value = "4295 Hope Aged Care Sunshine West"
$ws.Range("B9").Value = 30
$ws.Range("A10").Value = "44444 Nar Nar Goon Primary School Nar NarGoon"
$ws.Range("B10").Value = 12
$ws.Range("A11").Value = "44666 Gardenvale Primary School Senior SchoolCampus Brighton East"
$ws.Range("B11").Value = 18
$ws.Range("A12").Value = "44811 Dandenong North Primary SchoolDandenong"
$ws.Range("B12").Value = 14
$ws.Range("A13").Value = "44950 Templestowe Valley Primary SchoolTemplestowe Lower"
$ws.Range("B13").Value = 67
$ws.Range("A14").Value = "44979 Campbellfield Heights Primary SchoolCampbellfield"
$ws.Range("B14").Value = 11
$ws.Range("A15").Value = "44982 Diamond Creek East Primary SchoolDiamond Creek"
$ws.Range("B15").Value = 17
$ws.Range("A16").Value = "45248 Brookside P-9 College Caroline Springs"
$ws.Range("B16").Value = 30
$ws.Range("A17").Value = "45257 Roxburgh Rise Primary School RoxburghPark"
$ws.Range("B17").Value = 15
$ws.Range("A18").Value = "45267 Epping Views Primary School Epping"
$ws.Range("B18").Value = 15
$ws.Range("A19").Value = "45315 Red Hill Consolidated School Red Hill"
$ws.Range("B19").Value = 11
$ws.Range("A20").Value = "45585 Mount Ridley College Craigieburn"
$ws.Range("B20").Value = 12
$ws.Range("A21").Value = "45648 St Brendans Primary School Shepparton"
$ws.Range("B21").Value = 24
$ws.Range("A22").Value = "4574 Village Glen Aged Care ResidencesMornington"
$ws.Range("B22").Value = 11
$ws.Range("A23").Value = "45755 St Patricks Catholic Parish PrimarySchool Mentone"
$ws.Range("B23").Value = 18
$ws.Range("A24").Value = "45797 St John's Primary School Footscray"
$ws.Range("B24").Value = 13
$ws.Range("A25").Value = "45846 St Mary's School Mooroopna"
$ws.Range("B25").Value = 19
$ws.Range("A26").Value = "45903 St Peter Chanel Deer Park"
$ws.Range("B26").Value = 14
$ws.Range("A27").Value = "45950 St Luke's Primary School Lalor"
$ws.Range("B27").Value = 20
$ws.Range("A28").Value = "46001 Good Shepherd Parish School WheelersHill"
$ws.Range("B28").Value = 15
$ws.Range("A29").Value = "46052 St. Francis of Assisi Primary School MillPark"
$ws.Range("B29").Value = 26
$ws.Range("A30").Value = "46105 Christ the Priest Primary School CarolineSprings"
$ws.Range("B30").Value = 44
$ws.Range("A31").Value = "46115 St Luke's Catholic Primary SchoolShepparton North"
$ws.Range("B31").Value = 11
$ws.Range("A32").Value = "46117 Marymede Catholic College SouthMorang"
$ws.Range("B32").Value = 11
$ws.Range("A33").Value = "46125 Our Lady of the Southern Cross PrimarySchool Manor Lakes"
$ws.Range("B33").Value = 12
$ws.Range("A34").Value = "46221 Bialik College Hawthorn"
$ws.Range("B34").Value = 12
$ws.Range("A35").Value = "46239 Gilson College Taylors Hill"
$ws.Range("B35").Value = 10
$ws.Range("A36").Value = "50279 Dallas Brooks Community Primary SchoolDallas"
$ws.Range("B36").Value = 10
$ws.Range("A37").Value = "52476 Harvest Home Primary School Epping"
$ws.Range("B37").Value = 10
$ws.Range("A38").Value = "Covenant College Bell Post Hill"
$ws.Range("B38").Value = 15
$ws.Range("A39").Value = "House Party 27 November Private ResidenceBrunswick West"
$ws.Range("B39").Value = 26
$ws.Range("A40").Value = "Springside Primary School Caroline Springs Nov"
$ws.Range("B40").Value = 44
$ws.Range("A41").Value = "The Village Early Learning Centre Sandringham"
$ws.Range("B41").Value = 16
$ws.Range("A42").Value = "Torquay Hotel Torquay"
$ws.Range("B42").Value = 20
